$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new text values for row 3 (October 2nd entry)
$ws.Range("B3").Value = "Charlatan"
$ws.Range("E3").Value = "Affinité animale"
$ws.Range("H3").Value = "Costaud"
$ws.Range("I3").Value = "Armes à deux mains"
$ws.Range("J3").Value = "Éclaireur"

# Apply matching cell styles (Good/Neutral/Bad) by copying formatting
# from existing cells in row 2 that already carry the desired style,
# so the existing style indices are reused instead of new ones created.

# Neutral style -> B3
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# Good style -> E3, H3, I3, J3
$ws.Range("C2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)

# Bad style -> D3, K3, L3, M3
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("M3").PasteSpecial(-4122)

$excel.CutCopyMode = 0
